$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Set the value of the "Name" row (B4) which was previously empty
$ws.Range("B4").Value = "NatcycleformVs"

# Update the "Date" row (B8) to the new generation timestamp
$ws.Range("B8").Value = "2025-07-18T06:40:38+00:00"
